# Updates odds values in Sheet1 to match the 2024-10-14 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 changes
$ws.Range("G5").Value = 1.95
$ws.Range("I5").Value = 3.8
$ws.Range("L5").Value = 4.75
$ws.Range("N5").Value = 7.5
$ws.Range("AG5").Value = 9
$ws.Range("AH5").Value = 19

# Row 9 changes
$ws.Range("G9").Value = 2.9
$ws.Range("H9").Value = 3.3
$ws.Range("I9").Value = 2.15
$ws.Range("J9").Value = 3.6
$ws.Range("L9").Value = 2.88
$ws.Range("O9").Value = 1.25
$ws.Range("P9").Value = 3.75
$ws.Range("Q9").Value = 1.9
$ws.Range("R9").Value = 1.95
$ws.Range("Z9").Value = 34
$ws.Range("AH9").Value = 11
$ws.Range("AI9").Value = 9
$ws.Range("AJ9").Value = 21
$ws.Range("AK9").Value = 17
$ws.Range("AO9").Value = 17
$ws.Range("AW9").Value = 4.33
$ws.Range("AX9").Value = 12
